# Scraper re-run: two new columns, "height" and "weight", are inserted
# right before the existing "fantasy points" column.
#
# Before:  A=idx  B=rec_yds  C=rec_td  D=fumbles  E=fantasy points
# After:   A=idx  B=rec_yds  C=rec_td  D=fumbles  E=height  F=weight  G=fantasy points
#
# Concretely (since "fantasy points" was the last/right-most column, col 5 = E):
#   - the old column E (header + 16 data rows) shifts two columns right, to G
#   - the vacated E and new F columns get the new "height" (6.25) / "weight" (250)
#     header + constant values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldLastCol = 5   # column E ("fantasy points"), before the edit
$newLastCol = $oldLastCol + 2   # column G, after the edit

$lastRow = $ws.UsedRange.Rows.Count

# Stash the old "fantasy points" header text + data (col E) before overwriting it;
# .Formula returns the exact literal text of the cell so nothing is lost/rounded.
$oldHeader = $ws.Cells.Item(1, $oldLastCol).Formula
$oldValues = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $oldValues[$r] = $ws.Cells.Item($r, $oldLastCol).Formula
}

# Move "fantasy points" out to its new column (G).
$ws.Cells.Item(1, $newLastCol).Value = $oldHeader
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newLastCol).Formula = $oldValues[$r]
}

# Fill in the new "height" / "weight" columns (E / F) with the scraped constants.
$ws.Cells.Item(1, $oldLastCol).Value = "height"
$ws.Cells.Item(1, $oldLastCol + 1).Value = "weight"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $oldLastCol).Value = 6.25
    $ws.Cells.Item($r, $oldLastCol + 1).Value = 250
}

# The two new header cells should carry the same bold/centered/bordered look as
# the rest of row 1 - copy that formatting over from the neighboring header cell.
$ws.Cells.Item(1, $oldLastCol - 1).Copy()
$ws.Range($ws.Cells.Item(1, $oldLastCol), $ws.Cells.Item(1, $newLastCol)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
